# Revert "Merge branch 'wrong-xlsform-col'"
#
# Restores the "survey" sheet's column-C header from "label" back to
# "message" (the column-C header on "choices" stays "label" - it is a
# separate shared-string entry), moves the active selection on "survey"
# from C1 to C4, and restores the "choices" sheet's first-page numbering.

$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")
$choices = $wb.Worksheets.Item("choices")

# "survey" sheet: column C header goes back to "message"
$survey.Range("C1").Value = "message"

# Move the active cell / selection on "survey" from C1 to C4
$survey.Activate()
$survey.Range("C4").Select()

# "choices" sheet: restore first-page numbering
$choices.PageSetup.FirstPageNumber = 1
